# Insert a new row at the top of the sheet; this shifts every existing row
# down by one (old row 1 headers -> row 2, old row 2..100 data -> row 3..101).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()

# Populate the new row 1 with the numeric column-index markers (0-12).
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12

# Give the new row 1 the same (bold / bordered / centered) header style that
# row 2 (the old header row) already carries.
$ws.Range("A2:M2").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)

# The old header row (now row 2) loses its special style in the target, and
# its "thread_size" / "material_surface" helper labels are cleared out.
$ws.Range("A2:M2").Style = "Normal"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
